# MGS-7128 text change for 2023 rent increase template (tenant's response)
#
# 1) Collapse the cached MERGEFIELD display text for oldRentAmount / newRentAmount
#    («<run/>oldRentAmount<run/>» split across 3 runs + spellcheck proofErr markers)
#    down to a single «oldRentAmount» / «newRentAmount» run, matching how Word
#    re-writes a field result after it has been "typed over"/refreshed.
# 2) Reword the tenant's response paragraph in Part 3: replace
#    "propose to apply to a Rent Officer for a determination of the <tab>open market rent. "
#    with
#    "propose to make a referral to a Rent Officer to confirm whether the proposed
#     rent increase is in line with the rent cap. "

$d = $word.ActiveDocument

$chevOpen  = [char]171   # «
$chevClose = [char]187   # »

# --- 1a. «oldRentAmount» -------------------------------------------------
$oldField = $chevOpen + "oldRentAmount" + $chevClose
$found = $d.Content.Find.Execute(
    $oldField, $true, $true, $false, $false, $false, $true, 1, $false,
    $oldField, 2)
if (-not $found) {
    throw "Could not find oldRentAmount field result text"
}

# --- 1b. «newRentAmount» -------------------------------------------------
$newField = $chevOpen + "newRentAmount" + $chevClose
$found = $d.Content.Find.Execute(
    $newField, $true, $true, $false, $false, $false, $true, 1, $false,
    $newField, 2)
if (-not $found) {
    throw "Could not find newRentAmount field result text"
}

# --- 2. Reword the Rent Officer referral sentence ------------------------
$oldSentence = "propose to apply to a Rent Officer for a determination of the ^topen market rent. "
$newSentence = "propose to make a referral to a Rent Officer to confirm whether the proposed rent increase is in line with the rent cap. "
$found = $d.Content.Find.Execute(
    $oldSentence, $true, $true, $false, $false, $false, $true, 1, $false,
    $newSentence, 2)
if (-not $found) {
    throw "Could not find the Rent Officer referral sentence"
}
